{"js": "// Regenerate the document the way the site generator does:\n//   1. Drop the three boiler-plate \"Home\" / \"Back to Home\" /\n//      \"Download Word Document\" hyperlink paragraphs that used to sit\n//      in front of the real content.\n//   2. Every table in the body gets an explicit 100% preferred width\n//      (tblW type=\"pct\" w=\"5000\") instead of the old type=\"auto\" w=\"0\",\n//      with the <w:tblW> element now ordered before <w:tblStyle/>.\n//\n// The Word JS API has no direct way to flip a table's preferred-width\n// type to \"percent\" (Table.width is read-only and there is no\n// preferredWidthType-like member on Word.Table), so the table-width\n// part is done with the standard work-around for that gap: round-trip\n// the table through OOXML (getOoxml/insertOoxml) and patch the\n// <w:tblW> element in the fragment before re-inserting it.\n\n// --- 1. remove the leading Home / Back to Home / Download links ----\nconst leadParagraphs = context.document.body.paragraphs;\nleadParagraphs.load(\"items/text\");\nawait context.sync();\n\nconst dropTexts = [\"Home\", \"\u2190 Back to Home\", \"Download Word Document\"];\nfor (const target of dropTexts) {\n  const p = leadParagraphs.items.find((pp) => pp.text === target);\n  if (p) {\n    p.delete();\n  }\n}\nawait context.sync();\n\n// --- 2. force every table to 100% preferred width -------------------\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < tables.items.length; i++) {\n  const table = tables.items[i];\n  const range = table.getRange(\"Whole\");\n  const ooxml = range.getOoxml();\n  await context.sync();\n\n  let xml = ooxml.value;\n\n  // Swap the old auto/0 width for an explicit 100% (pct/5000) width.\n  xml = xml.replace(\n    /<w:tblW\\s+[^/]*\\/>/,\n    '<w:tblW w:type=\"pct\" w:w=\"5000\"/>'\n  );\n\n  range.insertOoxml(xml, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Regenerate the document the way the site generator does:\n#   1. Drop the three boiler-plate \"Home\" / \"Back to Home\" /\n#      \"Download Word Document\" hyperlink paragraphs that used to sit\n#      in front of the real content.\n#   2. Every table in the body gets an explicit 100% preferred width\n#      (tblW type=\"pct\" w=\"5000\") instead of the old type=\"auto\" w=\"0\".\n\n$d = $word.ActiveDocument\n\n# --- 1. remove the leading Home / Back to Home / Download links -----\n$dropTexts = @(\"Home\", \"\u2190 Back to Home\", \"Download Word Document\")\nforeach ($target in $dropTexts) {\n    $p = $d.Paragraphs(1)\n    $text = $p.Range.Text.TrimEnd()\n    if ($text -eq $target) {\n        $p.Range.Delete()\n    }\n}\n\n# --- 2. force every table to 100% preferred width ---------------------\nfor ($i = 1; $i -le $d.Tables.Count; $i++) {\n    $t = $d.Tables($i)\n    $t.PreferredWidthType = 2   # wdPreferredWidthPercent\n    $t.PreferredWidth = 250     # 250 -> 5000/50ths-of-percent == 100%\n}\n"}
